$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "syntax" -> "syntaxe" in the comment for the "Creation de profil" entry (row 25)
$ws.Range("E25").Value2 = "Permet de savoir si un profil existe déjà, si la syntaxe du nom est correct ou si le serveur de base de données est allumé."

# Duplicate the "Logique gestion de profil" task (row 27) into the next working day (row 28),
# matching the date formatting used on the rows above it.
$ws.Range("B27").Copy() | Out-Null
$ws.Range("B28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B28").Value2 = 43158
$ws.Range("C28").Value2 = "Logique gestion de profil"
$ws.Range("D28").Value2 = "1h30"
$ws.Range("E28").Value2 = "En cours"

# Move the selection to the newly added row, like a user who just finished typing it in.
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("E28:G28").Select() | Out-Null
